$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: TOTAL (5-17 ans) - update values only ---
$ws.Range("C2").Value = 2307204
$ws.Range("D2").Value = 47.5
$ws.Range("E2").Value = 1094955
$ws.Range("F2").Value = 47.4
$ws.Range("G2").Value = 1094654
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 97326
$ws.Range("J2").Value = 0.9
$ws.Range("K2").Value = 20270
$ws.Range("L2").Value = 52.5
$ws.Range("M2").Value = 1212250

# --- Row 3: non_pdi -> hote ---
$ws.Range("A3").Value = "hote (5-17 y.o.)"
$ws.Range("B3").Value = "hote"
$ws.Range("C3").Value = 2077436
$ws.Range("D3").Value = 48.1
$ws.Range("E3").Value = 998326
$ws.Range("F3").Value = 46.9
$ws.Range("G3").Value = 975307
$ws.Range("H3").Value = 4.2
$ws.Range("I3").Value = 87856
$ws.Range("J3").Value = 0.8
$ws.Range("K3").Value = 15946
$ws.Range("L3").Value = 51.9
$ws.Range("M3").Value = 1079110

# --- Row 4: pdi -> idp_host ---
$ws.Range("A4").Value = "idp_host (5-17 y.o.)"
$ws.Range("B4").Value = "idp_host"
$ws.Range("C4").Value = 125059
$ws.Range("D4").Value = 46.6
$ws.Range("E4").Value = 58310
$ws.Range("F4").Value = 48.5
$ws.Range("G4").Value = 60712
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 3871
$ws.Range("J4").Value = 1.7
$ws.Range("K4").Value = 2166
$ws.Range("L4").Value = 53.4
$ws.Range("M4").Value = 66749

# --- Row 5 (NEW): retourne (5-17 y.o.) ---
$ws.Range("A5").Value = "retourne (5-17 y.o.)"
$ws.Range("B5").Value = "retourne"
$ws.Range("C5").Value = 77546
$ws.Range("D5").Value = 36.3
$ws.Range("E5").Value = 28124
$ws.Range("F5").Value = 55.8
$ws.Range("G5").Value = 43254
$ws.Range("H5").Value = 5.4
$ws.Range("I5").Value = 4200
$ws.Range("J5").Value = 2.5
$ws.Range("K5").Value = 1967
$ws.Range("L5").Value = 63.7
$ws.Range("M5").Value = 49422

# --- Row 6 (NEW): idp_site (5-17 y.o.) ---
$ws.Range("A6").Value = "idp_site (5-17 y.o.)"
$ws.Range("B6").Value = "idp_site"
$ws.Range("C6").Value = 27164
$ws.Range("D6").Value = 37.5
$ws.Range("E6").Value = 10195
$ws.Range("F6").Value = 56.6
$ws.Range("G6").Value = 15381
$ws.Range("H6").Value = 5.1
$ws.Range("I6").Value = 1399
$ws.Range("J6").Value = 0.7
$ws.Range("K6").Value = 190
$ws.Range("L6").Value = 62.5
$ws.Range("M6").Value = 16970
